$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") rows 2-96 were incorrectly recorded as 63;
# recover the dropped data by restoring the correct value of 263.
$ws.Range("B2:B96").Value = 263
